$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency prices / 1h volume changes (and two row re-sorts)
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "42.933.89"
Set-TextValue "E2" "  +0.29%  "
Set-TextValue "D3" "2.359.74"
Set-TextValue "E3" "  +2.03%  "
Set-TextValue "E4" "  -0.09%  "
Set-TextValue "D5" "302.32"
Set-TextValue "E5" "  +0.16%  "
Set-TextValue "D6" "95.40"
Set-TextValue "E6" "  +0.17%  "
Set-TextValue "D7" "0.504"
Set-TextValue "E7" "  -0.34%  "
Set-TextValue "E8" "  -0.02%  "
Set-TextValue "D9" "0.485"
Set-TextValue "E9" "  -1.30%  "
Set-TextValue "D10" "33.80"
Set-TextValue "E10" "  -1.17%  "
Set-TextValue "D11" "0.0783"
Set-TextValue "E11" "  +0.06%  "
Set-TextValue "E12" "  +3.01%  "
Set-TextValue "D13" "18.31"
Set-TextValue "E13" "  -3.18%  "
Set-TextValue "B14" "WrappedliquidstakedEther2.0"
Set-TextValue "C14" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D14" "2.728.57"
Set-TextValue "E14" "  +1.99%  "
Set-TextValue "B15" "Polkadot"
Set-TextValue "C15" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D15" "6.70"
Set-TextValue "E15" "  -0.48%  "
Set-TextValue "D16" "2.346.79"
Set-TextValue "E16" "  +1.76%  "
Set-TextValue "D17" "0.793"
Set-TextValue "E17" "  +0.91%  "
Set-TextValue "D18" "42.896.14"
Set-TextValue "E18" "  +0.34%  "
Set-TextValue "D19" "11.82"
Set-TextValue "E19" "  -2.80%  "
Set-TextValue "D20" "6.24"
Set-TextValue "E20" "  +1.77%  "
Set-TextValue "D21" "0.0₃0883"
Set-TextValue "D22" "67.84"
Set-TextValue "E22" "  +0.18%  "
Set-TextValue "D23" "235.03"
Set-TextValue "E23" "  -0.11%  "
Set-TextValue "D24" "2.18"
Set-TextValue "E24" "  -3.86%  "
Set-TextValue "E25" "  -0.06%  "
Set-TextValue "D26" "2.41"
Set-TextValue "E26" "  +0.16%  "
Set-TextValue "D27" "24.60"
Set-TextValue "E27" "  +1.24%  "
Set-TextValue "E28" "  +0.38%  "
Set-TextValue "D29" "9.22"
Set-TextValue "E29" "  +1.12%  "
Set-TextValue "D30" "31.45"
Set-TextValue "E30" "  -1.74%  "
Set-TextValue "E31" "  -0.07%  "
Set-TextValue "D32" "5.01"
Set-TextValue "E32" "  +0.25%  "
Set-TextValue "D33" "17.28"
Set-TextValue "E33" "  -2.42%  "
Set-TextValue "D34" "0.0716"
Set-TextValue "E34" "  +2.79%  "
Set-TextValue "E35" "  +3.51%  "
Set-TextValue "D36" "1.83"
Set-TextValue "E36" "  +3.13%  "
Set-TextValue "D37" "4.33"
Set-TextValue "E37" "  -2.77%  "
Set-TextValue "D38" "2.29"
Set-TextValue "E38" "  -1.77%  "
Set-TextValue "B39" "Monero"
Set-TextValue "C39" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D39" "121.88"
Set-TextValue "E39" "  -26.50%  "
Set-TextValue "B40" "LidoDAOToken"
Set-TextValue "C40" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D40" "2.76"
Set-TextValue "E40" "  +1.77%  "
Set-TextValue "E41" "  -0.85%  "
Set-TextValue "D42" "20.96"
Set-TextValue "E42" "  -0.10%  "
Set-TextValue "D43" "1.929.59"
Set-TextValue "E43" "  +0.18%  "
Set-TextValue "D44" "0.0278"
Set-TextValue "E44" "  -0.09%  "
Set-TextValue "E45" "  +2.43%  "
Set-TextValue "D46" "9.13"
Set-TextValue "E46" "  -9.89%  "
Set-TextValue "D47" "2.70"
Set-TextValue "E47" "  -1.60%  "
Set-TextValue "D48" "2.587.62"
Set-TextValue "E48" "  +1.73%  "
Set-TextValue "E49" "  +1.87%  "
Set-TextValue "B50" "TrustWalletToken"
Set-TextValue "C50" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D50" "1.14"
Set-TextValue "E50" "  +1.20%  "
Set-TextValue "D51" "51.51"
Set-TextValue "E51" "  -3.11%  "

Write-Output "Applied 99 cell updates"
